# Updates cryptos price/volume columns (D, E) per latest scrape run.
# Values that look numeric are force-written as Text (matching the
# original inline-string cell type) by temporarily applying a Text
# number format, then clearing formatting again so no residual style
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.362.64'
$ws.Range("D3").Value = '1.844.33'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.09'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6346'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2959'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.67'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.83%  '
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.982'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009893'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '29.380.94'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.50'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.35%  '
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9998'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.540'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  +237.28%  '
$ws.Range("E24").Value = '  +170.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.22'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1405'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.365'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05705'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.253'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.026'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.842'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.58%  '
$ws.Range("E35").Value = '  -1.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7153'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '1.249.09'
$ws.Range("E38").Value = '  +2.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.801'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01812'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.91%  '
$ws.Range("E41").Value = '  +267.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9016'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9995'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.78'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.064'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.147'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4017'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.702'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05739'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.25%  '
